# "change process CDM mapping"
# The observation_source_value row now maps to KTAS instead of the old
# LOINC code, and its source_concept_id is reset to the unmapped value 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "KTAS"
$ws.Range("B2").Value = 0

# Printer/page setup picked up on this save (portrait, paper size 9 = A4).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
